$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 2021 Form 1040 publish: line 27 ("Earned income credit (EIC)") is split
# into 27a/27b/27c, everything below shifts down by two rows, and the
# "total other payments" line picks up the new 27a + 28-31 formula.
# ---------------------------------------------------------------------------

# Insert two blank rows directly below the old line-27 row (row 43) so the
# rest of the block (old rows 44-56) shifts down to new rows 46-58.
$ws.Rows("44:45").Insert()

# Give the two new rows the same look as the line they were split from.
$ws.Range("A43:F43").Copy()
$ws.Range("A44:F45").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 43: was "27", now "27a" (EIC) -------------------------------------
$ws.Range("A43").Value = "27a"
$ws.Range("C43").Value = "27a"

# --- Row 44 (new): "27b" Nontaxable combat pay election ---------------------
$ws.Range("A44").Value = "27b"
$ws.Range("B44").Value = "Nontaxable combat pay election"
$ws.Range("C44").Value = "27b"

# --- Row 45 (new): "27c" Prior year (2019) earned income --------------------
$ws.Range("A45").Value = "27c"
$ws.Range("B45").Value = "Prior year (2019) earned income"
$ws.Range("C45").Value = "27c"

# --- Row 49 (old row 47): Schedule 3 reference updated 13 -> 15 -------------
$ws.Range("B49").Value = "Schedule 3, line 15"

# --- Row 50 (old row 48): total other payments & credits, new formula -------
$ws.Range("B50").Value = "Add lines 27a and 28 through 31. These are your total other payments and refundable credits"
$ws.Range("F50").Formula = "=D43+SUM(D46:D49)"

# ---------------------------------------------------------------------------
# View / selection tidy-up to match the published workbook.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.TopLeftCell = $ws.Range("A1")
$ws.Range("F50").Select()
